$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.277433666666667
$ws.Range("H2").Value = 3.832301
$ws.Range("I2").Value = 0.01786062203930835
$ws.Range("J2").Value = 0.01786062203930835
$ws.Range("M2").Value = 41.83357366666667
$ws.Range("N2").Value = 125.500721
$ws.Range("O2").Value = 0.2773195847425811
$ws.Range("P2").Value = 0.2773195847425811
$ws.Range("Q2").Value = 53.43961539878011
$ws.Range("R2").Value = 480.956538589021
$ws.Range("S2").Value = 0.004953100287185186
$ws.Range("T2").Value = 0.004953100287185183
$ws.Range("G3").Value = 1.277433666666667
$ws.Range("H3").Value = 3.832301
$ws.Range("I3").Value = 0.01786062203930835
$ws.Range("J3").Value = 0.01786062203930835
$ws.Range("O3").Value = 0.4239803668761465
$ws.Range("P3").Value = 0.4239803668761465
$ws.Range("Q3").Value = 81.70121761694688
$ws.Range("R3").Value = 735.3109585525219
$ws.Range("S3").Value = 0.007572553084862144
$ws.Range("T3").Value = 0.007572553084862143
$ws.Range("G4").Value = 1.277433666666667
$ws.Range("H4").Value = 3.832301
$ws.Range("I4").Value = 0.01786062203930835
$ws.Range("J4").Value = 0.01786062203930835
$ws.Range("M4").Value = 14.18032733333333
$ws.Range("N4").Value = 42.540982
$ws.Range("O4").Value = 0.09400302539123752
$ws.Range("P4").Value = 0.09400302539123752
$ws.Range("Q4").Value = 18.11442753995355
$ws.Range("R4").Value = 163.029847859582
$ws.Range("S4").Value = 0.0016789525070644
$ws.Range("T4").Value = 0.001678952507064399
$ws.Range("G5").Value = 1.277433666666667
$ws.Range("H5").Value = 3.832301
$ws.Range("I5").Value = 0.01786062203930835
$ws.Range("J5").Value = 0.01786062203930835
$ws.Range("M5").Value = 30.87848266666667
$ws.Range("N5").Value = 92.635448
$ws.Range("O5").Value = 0.204697022990035
$ws.Range("P5").Value = 0.2046970229900349
$ws.Range("Q5").Value = 39.44521333398311
$ws.Range("R5").Value = 355.006920005848
$ws.Range("S5").Value = 0.003656016160196628
$ws.Range("T5").Value = 0.003656016160196627
$ws.Range("G6").Value = 45.44725166666667
$ws.Range("I6").Value = 0.6354272679079697
$ws.Range("J6").Value = 0.6354272679079697
$ws.Range("M6").Value = 41.83357366666667
$ws.Range("N6").Value = 125.500721
$ws.Range("O6").Value = 0.2773195847425811
$ws.Range("P6").Value = 0.2773195847425811
$ws.Range("Q6").Value = 1901.220950545039
$ws.Range("R6").Value = 17110.98855490535
$ws.Range("S6").Value = 0.176216426070351
$ws.Range("T6").Value = 0.176216426070351
$ws.Range("G7").Value = 45.44725166666667
$ws.Range("I7").Value = 0.6354272679079697
$ws.Range("J7").Value = 0.6354272679079697
$ws.Range("O7").Value = 0.4239803668761465
$ws.Range("P7").Value = 0.4239803668761465
$ws.Range("R7").Value = 26160.15458070311
$ws.Range("S7").Value = 0.2694086861707284
$ws.Range("T7").Value = 0.2694086861707284
$ws.Range("G8").Value = 45.44725166666667
$ws.Range("I8").Value = 0.6354272679079697
$ws.Range("J8").Value = 0.6354272679079697
$ws.Range("M8").Value = 14.18032733333333
$ws.Range("N8").Value = 42.540982
$ws.Range("O8").Value = 0.09400302539123752
$ws.Range("P8").Value = 0.09400302539123752
$ws.Range("Q8").Value = 644.4569050337121
$ws.Range("R8").Value = 5800.11214530341
$ws.Range("S8").Value = 0.05973208559943756
$ws.Range("T8").Value = 0.05973208559943756
$ws.Range("G9").Value = 45.44725166666667
$ws.Range("I9").Value = 0.6354272679079697
$ws.Range("J9").Value = 0.6354272679079697
$ws.Range("M9").Value = 30.87848266666667
$ws.Range("N9").Value = 92.635448
$ws.Range("O9").Value = 0.204697022990035
$ws.Range("P9").Value = 0.2046970229900349
$ws.Range("Q9").Value = 1403.342172836804
$ws.Range("R9").Value = 12630.07955553124
$ws.Range("S9").Value = 0.1300700700674528
$ws.Range("T9").Value = 0.1300700700674528
$ws.Range("G10").Value = 23.96074166666667
$ws.Range("H10").Value = 71.88222500000001
$ws.Range("I10").Value = 0.3350105464235513
$ws.Range("J10").Value = 0.3350105464235513
$ws.Range("M10").Value = 41.83357366666667
$ws.Range("N10").Value = 125.500721
$ws.Range("O10").Value = 0.2773195847425811
$ws.Range("P10").Value = 0.2773195847425811
$ws.Range("Q10").Value = 1002.36345162047
$ws.Range("R10").Value = 9021.271064584225
$ws.Range("S10").Value = 0.09290498561856445
$ws.Range("T10").Value = 0.09290498561856442
$ws.Range("G11").Value = 23.96074166666667
$ws.Range("H11").Value = 71.88222500000001
$ws.Range("I11").Value = 0.3350105464235513
$ws.Range("J11").Value = 0.3350105464235513
$ws.Range("O11").Value = 0.4239803668761465
$ws.Range("P11").Value = 0.4239803668761465
$ws.Range("Q11").Value = 1532.464518709605
$ws.Range("R11").Value = 13792.18066838645
$ws.Range("S11").Value = 0.1420378943800356
$ws.Range("T11").Value = 0.1420378943800356
$ws.Range("G12").Value = 23.96074166666667
$ws.Range("H12").Value = 71.88222500000001
$ws.Range("I12").Value = 0.3350105464235513
$ws.Range("J12").Value = 0.3350105464235513
$ws.Range("M12").Value = 14.18032733333333
$ws.Range("N12").Value = 42.540982
$ws.Range("O12").Value = 0.09400302539123752
$ws.Range("P12").Value = 0.09400302539123752
$ws.Range("Q12").Value = 339.7711599827722
$ws.Range("R12").Value = 3057.94043984495
$ws.Range("S12").Value = 0.03149200490178545
$ws.Range("T12").Value = 0.03149200490178545
$ws.Range("G13").Value = 23.96074166666667
$ws.Range("H13").Value = 71.88222500000001
$ws.Range("I13").Value = 0.3350105464235513
$ws.Range("J13").Value = 0.3350105464235513
$ws.Range("M13").Value = 30.87848266666667
$ws.Range("N13").Value = 92.635448
$ws.Range("O13").Value = 0.204697022990035
$ws.Range("P13").Value = 0.2046970229900349
$ws.Range("Q13").Value = 739.8713462346444
$ws.Range("R13").Value = 6658.8421161118
$ws.Range("S13").Value = 0.06857566152316587
$ws.Range("T13").Value = 0.06857566152316584
$ws.Range("G14").Value = 0.8369233333333334
$ws.Range("H14").Value = 2.51077
$ws.Range("I14").Value = 0.01170156362917063
$ws.Range("J14").Value = 0.01170156362917063
$ws.Range("M14").Value = 41.83357366666667
$ws.Range("N14").Value = 125.500721
$ws.Range("O14").Value = 0.2773195847425811
$ws.Range("P14").Value = 0.2773195847425811
$ws.Range("Q14").Value = 35.01149391835222
$ws.Range("R14").Value = 315.10344526517
$ws.Range("S14").Value = 0.00324507276648049
$ws.Range("T14").Value = 0.003245072766480489
$ws.Range("G15").Value = 0.8369233333333334
$ws.Range("H15").Value = 2.51077
$ws.Range("I15").Value = 0.01170156362917063
$ws.Range("J15").Value = 0.01170156362917063
$ws.Range("O15").Value = 0.4239803668761465
$ws.Range("P15").Value = 0.4239803668761465
$ws.Range("Q15").Value = 53.52736284443777
$ws.Range("R15").Value = 481.7462655999399
$ws.Range("S15").Value = 0.004961233240520337
$ws.Range("T15").Value = 0.004961233240520336
$ws.Range("G16").Value = 0.8369233333333334
$ws.Range("H16").Value = 2.51077
$ws.Range("I16").Value = 0.01170156362917063
$ws.Range("J16").Value = 0.01170156362917063
$ws.Range("M16").Value = 14.18032733333333
$ws.Range("N16").Value = 42.540982
$ws.Range("O16").Value = 0.09400302539123752
$ws.Range("P16").Value = 0.09400302539123752
$ws.Range("Q16").Value = 11.86784681957111
$ws.Range("R16").Value = 106.81062137614
$ws.Range("S16").Value = 0.001099982382950108
$ws.Range("T16").Value = 0.001099982382950108
$ws.Range("G17").Value = 0.8369233333333334
$ws.Range("H17").Value = 2.51077
$ws.Range("I17").Value = 0.01170156362917063
$ws.Range("J17").Value = 0.01170156362917063
$ws.Range("M17").Value = 30.87848266666667
$ws.Range("N17").Value = 92.635448
$ws.Range("O17").Value = 0.204697022990035
$ws.Range("P17").Value = 0.2046970229900349
$ws.Range("Q17").Value = 25.84292264166222
$ws.Range("R17").Value = 232.58630377496
$ws.Range("S17").Value = 0.002395275239219698
$ws.Range("T17").Value = 0.002395275239219698
